$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume/1h-change (E) columns for rows with new data.
# D-column cells are stored as TEXT in the source data (e.g. "65.805.26" uses
# dot-grouping, not a decimal point); where the new value would otherwise be
# auto-parsed as a number by Excel, force the cell to Text format first so it
# round-trips as a string exactly like the original.
$ws.Range("D2").Value = "65.805.26"
$ws.Range("E2").Value = "  +6.18%  "
$ws.Range("D3").Value = "3.006.73"
$ws.Range("E3").Value = "  +3.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.35"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.61"
$ws.Range("E6").Value = "  +12.42%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +3.45%  "
$ws.Range("D9").Value = "3.004.33"
$ws.Range("E9").Value = "  +3.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("E10").Value = "  -4.99%  "
$ws.Range("E11").Value = "  +6.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  +7.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  +8.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.60"
$ws.Range("E14").Value = "  +7.67%  "
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "65.822.57"
$ws.Range("E16").Value = "  +6.28%  "
$ws.Range("D17").Value = "3.506.98"
$ws.Range("E17").Value = "  +3.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.97"
$ws.Range("E18").Value = "  +7.36%  "
$ws.Range("D19").Value = "3.015.33"
$ws.Range("E19").Value = "  +3.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "456.82"
$ws.Range("E20").Value = "  +6.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.94"
$ws.Range("E21").Value = "  +8.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.688"
$ws.Range("E22").Value = "  +5.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.36"
$ws.Range("E23").Value = "  +7.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.29"
$ws.Range("E24").Value = "  +4.20%  "
$ws.Range("E25").Value = "  +12.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.37"
$ws.Range("E26").Value = "  +3.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.66"
$ws.Range("E27").Value = "  +4.96%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  +16.67%  "
$ws.Range("E30").Value = "  +15.87%  "
$ws.Range("E31").Value = "  -7.26%  "
$ws.Range("E32").Value = "  +3.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.95"
$ws.Range("E33").Value = "  +5.21%  "
$ws.Range("E34").Value = "  +3.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +3.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.80"
$ws.Range("E37").Value = "  +7.75%  "
$ws.Range("E38").Value = "  +13.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.72"
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.95"
$ws.Range("E40").Value = "  +1.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.310"
$ws.Range("E41").Value = "  +15.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.122"
$ws.Range("E42").Value = "  +6.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.83"
$ws.Range("E43").Value = "  +6.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.44"
$ws.Range("E44").Value = "  +3.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "391.86"
$ws.Range("E45").Value = "  +13.16%  "

# Rows 46/47: VeChain and Maker swapped order, with refreshed price/volume figures
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.789.18"
$ws.Range("E46").Value = "  +3.17%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0355"
$ws.Range("E47").Value = "  +5.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.96"
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.80"
$ws.Range("E50").Value = "  +10.22%  "
$ws.Range("E51").Value = "  +4.10%  "
